$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Đại vận xấu, sẽ phải cố gắng lắm nhưng cũng để trả nợ đời. Âu cũng phải cố gắng vì làm gì có nợ nào mà không phải trả."
$ws.Range("B2").Value = "Đại vận sinh xuất gặp nhiều lao đao, nhưng có nghị lực vẫn gặt quả ngọt."
$ws.Range("B3").Value = "Đại vận đại cát sẽ gặp nhiều việc tốt, nhiều cơ hội, thuận lợi với bản thân."
$ws.Range("B4").Value = "Đại vận rất xấu, sẽ gặp nhiều khó khăn và buồn lo. Nhưng sẽ có nhiều bài học sâu sắc."
$ws.Range("B5").Value = "Đại vận đẹp sẽ gặp nhiều việc toại ý."

$ws.Range("B4").Select()
